$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.067.23"
$ws.Range("E2").Value = "  -0.54%  "
$ws.Range("D3").Value = "2.548.21"
$ws.Range("E3").Value = "  -0.15%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'585.01"
$ws.Range("E5").Value = "  +2.20%  "
$ws.Range("D6").Value = "'147.32"
$ws.Range("E6").Value = "  -2.50%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").Value = "'0.584"
$ws.Range("E8").Value = "  -1.07%  "
$ws.Range("E9").Value = "  -0.65%  "
$ws.Range("E10").Value = "  -3.38%  "
$ws.Range("E11").Value = "  -0.23%  "
$ws.Range("E12").Value = "  -1.21%  "
$ws.Range("D13").Value = "'27.47"
$ws.Range("E13").Value = "  -3.61%  "
$ws.Range("D14").Value = "2.996.37"
$ws.Range("E14").Value = "  -0.32%  "
$ws.Range("D15").Value = "62.944.30"
$ws.Range("E15").Value = "  -0.62%  "
$ws.Range("E16").Value = "  -0.88%  "
$ws.Range("D17").Value = "2.545.65"
$ws.Range("E17").Value = "  +0.73%  "
$ws.Range("D18").Value = "'11.34"
$ws.Range("E18").Value = "  -2.91%  "
$ws.Range("D19").Value = "'337.43"
$ws.Range("E19").Value = "  -0.92%  "
$ws.Range("E20").Value = "  -1.18%  "
$ws.Range("E21").Value = "  -1.40%  "
$ws.Range("E22").Value = "  +0.26%  "
$ws.Range("D23").Value = "'65.93"
$ws.Range("E23").Value = "  -0.32%  "
$ws.Range("E24").Value = "  -0.72%  "
$ws.Range("E25").Value = "  +0.63%  "
$ws.Range("B26").Value = "Binance-PegBSC-USD"
$ws.Range("C26").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D26").Value = "'1.01"
$ws.Range("E26").Value = "  +1.08%  "
$ws.Range("B27").Value = "SuiNetwork"
$ws.Range("C27").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D27").Value = "'1.49"
$ws.Range("E27").Value = "  +0.99%  "
$ws.Range("D28").Value = "'8.39"
$ws.Range("E28").Value = "  -1.14%  "
$ws.Range("D29").Value = "'7.72"
$ws.Range("E29").Value = "  +9.06%  "
$ws.Range("D30").Value = "'1.97"
$ws.Range("E30").Value = "  +5.43%  "
$ws.Range("D31").Value = "0.0₃0815"
$ws.Range("E31").Value = "  -2.93%  "
$ws.Range("D32").Value = "'178.17"
$ws.Range("E32").Value = "  +0.90%  "
$ws.Range("B33").Value = "Bittensor"
$ws.Range("C33").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D33").Value = "'418.25"
$ws.Range("E33").Value = "  -0.90%  "
$ws.Range("B34").Value = "ImmutableX"
$ws.Range("C34").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D34").Value = "'1.54"
$ws.Range("E34").Value = "  -0.99%  "
$ws.Range("D35").Value = "'19.19"
$ws.Range("E35").Value = "  -0.12%  "
$ws.Range("E36").Value = "  -1.70%  "
$ws.Range("E37").Value = "  +0.03%  "
$ws.Range("D38").Value = "'4.35"
$ws.Range("E38").Value = "  -2.51%  "
$ws.Range("D39").Value = "'1.75"
$ws.Range("E39").Value = "  -1.23%  "
$ws.Range("E40").Value = "  +0.07%  "
$ws.Range("D41").Value = "'39.62"
$ws.Range("E41").Value = "  -0.97%  "
$ws.Range("D42").Value = "'150.72"
$ws.Range("E42").Value = "  -3.04%  "
$ws.Range("D43").Value = "'3.78"
$ws.Range("E43").Value = "  -0.75%  "
$ws.Range("E44").Value = "  -1.68%  "
$ws.Range("D45").Value = "'0.0541"
$ws.Range("E45").Value = "  +1.61%  "
$ws.Range("B46").Value = "Stellar"
$ws.Range("C46").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D46").Value = "'0.0974"
$ws.Range("E46").Value = "  +0.58%  "
$ws.Range("B47").Value = "Mantle"
$ws.Range("C47").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D47").Value = "'0.602"
$ws.Range("E47").Value = "  -1.51%  "
$ws.Range("D48").Value = "'0.0238"
$ws.Range("E48").Value = "  -1.43%  "
$ws.Range("D49").Value = "'18.29"
$ws.Range("E49").Value = "  -2.02%  "
$ws.Range("D50").Value = "'1.72"
$ws.Range("E50").Value = "  -6.34%  "
$ws.Range("E51").Value = "  -0.27%  "
